$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Update the query text held in the shared strings (B2, C2:C4, B3, B4)
# ---------------------------------------------------------------------------

$statQuery = @'
CALL{
    MATCH (p:participant)-->(s:study)
    OPTIONAL MATCH (samp:sample)-->(p)
    OPTIONAL MATCH (samp)<--(f:file)
    OPTIONAL MATCH (f)<--(g:genomic_info)
    OPTIONAL MATCH (p)<--(diag:diagnosis)
    WITH s, p, samp, f, g, diag, apoc.coll.flatten(COLLECT (apoc.text.split(f.experimental_strategy_and_data_subtypes,"[;,]\s{0,1}")), true) as es
    WHERE "Archer Fusion" IN es
    RETURN 
        count(distinct p) AS num_participants
}
WITH num_participants
CALL {
    MATCH (samp:sample)-->(p:participant)-->(s:study)
    OPTIONAL MATCH (samp)<--(f:file)
    OPTIONAL MATCH (p)<--(diag:diagnosis)
    OPTIONAL MATCH (f)<--(g:genomic_info)
    OPTIONAL MATCH (p)<--(diag:diagnosis)
    WITH s, p, samp, f, g, diag, apoc.coll.flatten(COLLECT (apoc.text.split(f.experimental_strategy_and_data_subtypes,"[;,]\s{0,1}")), true) as es
    WHERE "Archer Fusion" IN es
    RETURN 
        count(distinct samp) AS num_samples
}
WITH num_participants, num_samples
CALL {
    MATCH (f:file)-->(s:study)
    OPTIONAL MATCH (f)<--(g:genomic_info)
    OPTIONAL MATCH (samp:sample)<--(f)
    OPTIONAL MATCH (p:participant)<--(samp)
    OPTIONAL MATCH (p)<--(diag:diagnosis)
    WITH s, p, samp, f, g, diag, apoc.coll.flatten(COLLECT (apoc.text.split(f.experimental_strategy_and_data_subtypes,"[;,]\s{0,1}")), true) as es
    WHERE "Archer Fusion" IN es
    RETURN 
        count(distinct s) AS num_studies,
        count(distinct f) AS num_files
}
RETURN 
    num_studies AS Studies,
    num_participants AS Participants,
    num_samples AS Samples,
    num_files AS `Files`
'@

$participantsQuery = @'
MATCH (p:participant)-->(s:study)
OPTIONAL MATCH (samp:sample)-->(p)
OPTIONAL MATCH (samp)<--(f:file)
WITH p, samp, apoc.coll.flatten(COLLECT (apoc.text.split(f.experimental_strategy_and_data_subtypes,"[;,]\s{0,1}")), true) as es
WHERE "Archer Fusion" IN es
WITH p
OPTIONAL MATCH (p)-->(s:study)
OPTIONAL MATCH (samp:sample)-->(p)
WITH s, p, apoc.coll.sort(collect(distinct coalesce(samp.sample_id, "Not specified in data"))) as samp
RETURN 
coalesce(p.participant_id,'') as `Participant ID`,
coalesce(s.study_name, '') as `Study Name`,
coalesce(s.phs_accession,'') as `Accession`,
coalesce(p.gender,'') as `Gender`,
coalesce(apoc.text.join(samp, ','), '') as `Samples`
ORDER BY p.participant_id LIMIT 100
'@

$samplesQuery = @'
MATCH (samp:sample)-->(p:participant)-->(s:study)
OPTIONAL MATCH (samp)<--(f:file)
OPTIONAL MATCH (f)<--(g:genomic_info)
OPTIONAL MATCH (p)<--(diag:diagnosis)
WITH s, p, samp, f, g, diag, apoc.coll.flatten(COLLECT (apoc.text.split(f.experimental_strategy_and_data_subtypes,"[;,]\s{0,1}")), true) as es
WHERE "Archer Fusion" IN es
WITH DISTINCT s, p, samp
RETURN
    coalesce(samp.sample_id, '') as `Sample ID`,
    coalesce(p.participant_id,'') as `Participant ID`,
    coalesce(s.study_name, '') as `Study Name`,
    coalesce(s.phs_accession,'') as `Accession`,
    coalesce(samp.sample_tumor_status,'') as `Tumor`,
    coalesce(samp.sample_type,'') as `Analyte Type`
ORDER BY samp.sample_id LIMIT 100
'@

$filesQuery = @'
MATCH (f:file)-->(s:study)
OPTIONAL MATCH (samp:sample)<--(f)
OPTIONAL MATCH (samp)-->(p:participant)
OPTIONAL MATCH (f)<--(g:genomic_info)
OPTIONAL MATCH (p)<--(diag:diagnosis)
WITH s, p, samp, f, g, diag, apoc.coll.flatten(COLLECT (apoc.text.split(f.experimental_strategy_and_data_subtypes,"[;,]\s{0,1}")), true) as es
WHERE "Archer Fusion" IN es
WITH DISTINCT f, s, p, samp
RETURN
    coalesce(f.file_name, '') as `File Name`,
    coalesce(s.study_name,'') as `Study Name`,
    coalesce(s.phs_accession,'') as `Accession`,
    coalesce(p.participant_id, '') as `Participant ID`,
    coalesce(samp.sample_id, '') as `Sample ID`,
    coalesce(f.file_type, '') as `File Type`
ORDER BY f.file_name LIMIT 100
'@

$ws.Range("C2:C4").Value = $statQuery
$ws.Range("B2").Value = $participantsQuery
$ws.Range("B3").Value = $samplesQuery
$ws.Range("B4").Value = $filesQuery

# ---------------------------------------------------------------------------
# 2. Increase the shared font size used by the query cells from 12 to 14 and
#    apply that font (without wrap) to every other populated cell.
# ---------------------------------------------------------------------------

$ws.Range("B2:C4").Font.Size = 14

$ws.Range("A1:E1").Font.Size = 14
$ws.Range("A1:E1").WrapText = $false

$ws.Range("A2").Font.Size = 14
$ws.Range("A2").WrapText = $false
$ws.Range("D2:E2").Font.Size = 14
$ws.Range("D2:E2").WrapText = $false

$ws.Range("A3").Font.Size = 14
$ws.Range("A3").WrapText = $false
$ws.Range("D3:E3").Font.Size = 14
$ws.Range("D3:E3").WrapText = $false

$ws.Range("A4").Font.Size = 14
$ws.Range("A4").WrapText = $false
$ws.Range("D4:E4").Font.Size = 14
$ws.Range("D4:E4").WrapText = $false

$ws.Range("B5:C5").Font.Size = 14
$ws.Range("B5:C5").WrapText = $true
$ws.Range("C6").Font.Size = 14
$ws.Range("C6").WrapText = $true

# ---------------------------------------------------------------------------
# 3. Row heights / column widths / view settings
# ---------------------------------------------------------------------------

$ws.Rows.Item(2).RowHeight = 409.5
$ws.Rows.Item(3).RowHeight = 409.5
$ws.Rows.Item(4).RowHeight = 409.5

$ws.Columns.Item(1).ColumnWidth = 12.28515625
$ws.Columns.Item(2).ColumnWidth = 75.7109375
$ws.Columns.Item(3).ColumnWidth = 75.7109375
$ws.Columns.Item(4).ColumnWidth = 70.28515625
$ws.Columns.Item(5).ColumnWidth = 68.7109375

$ws.Range("E4").Select()
$excel.ActiveWindow.ScrollRow = 2
$excel.ActiveWindow.ScrollColumn = 1
